# Auto-generated variable block: content strings used by the edit
$WebExcel = @'
WebExcel
'@
$query = @'
query
'@
$dbExcel = @'
dbExcel
'@
$StatQuery = @'
StatQuery
'@
$TabName = @'
TabName
'@
$CasesTab = @'
CasesTab
'@
$SamplesTab = @'
SamplesTab
'@
$FilesTab = @'
FilesTab
'@
$StudyFilesTab = @'
StudyFilesTab
'@
$Neo4jData = @'
TC09_Canine_Filter_SamplePatho-PulmoCarcinoma_Neo4jData.xlsx
'@
$WebData = @'
TC09_Canine_Filter_SamplePatho-PulmoCarcinoma_WebData.xlsx
'@
$cases_query = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
MATCH (samp:sample)-->(c)
WHERE  samp.specific_sample_pathology in ['Pulmonary Carcinoma']
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@
$stat_query_broken = @'
MATCH (p:program)<--(s:study)<--(c)
MATCH (cf)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Pulmonary Carcinoma]
MATCH (cf:file)-[*]->(c:case)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
   count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$stat_query_fixed = @'
MATCH (p:program)<--(s:study)<--(c)
MATCH (cf)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Pulmonary Carcinoma']
MATCH (cf:file)-[*]->(c:case)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$files_query = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
MATCH (f)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Pulmonary Carcinoma']
 MATCH (f)-[*]->(samp:sample)
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 200
'@
$samples_query = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE  samp.specific_sample_pathology IN ['Pulmonary Carcinoma']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 200
'@
$study_files_query = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (c)<--(demo:demographic)
MATCH (samp:sample)-->(c)
WHERE samp.specific_sample_pathology IN ['Pulmonary Carcinoma']
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (CasesTab): new Cases query + (still-broken) stat query ---
$ws.Range("B2").Value = $cases_query
$ws.Range("C2").Value = $stat_query_broken

# --- Row 3 (SamplesTab): reworded Samples query + fixed stat query ---
$ws.Range("B3").Value = $samples_query
$ws.Range("C3").Value = $stat_query_fixed

# --- Row 4 (FilesTab): rewritten Files query + fixed stat query ---
$ws.Range("B4").Value = $files_query
$ws.Range("C4").Value = $stat_query_fixed

# --- Row 5 (new StudyFilesTab) ---
$ws.Range("A5").Value = $StudyFilesTab
$ws.Range("B5").Value = $study_files_query
$ws.Range("C5").Value = $stat_query_fixed
$ws.Range("D5").Value = $Neo4jData
$ws.Range("E5").Value = $WebData

# --- Formatting: wrap text on the long query/tool columns ---
$ws.Range("B2:C5").WrapText = $true

# C column (StatQuery) gets vertically centered wrapped text
$ws.Range("C2:C4").VerticalAlignment = -4108

# B3, B4, B5 (and C5) use the larger 18pt font introduced for this edit
$ws.Range("B3:B5").Font.Size = 18
$ws.Range("C5").Font.Size = 18

# --- Row heights to fit the newly-expanded content ---
$ws.Rows(2).RowHeight = 304.5
$ws.Rows(3).RowHeight = 409.5
$ws.Rows(4).RowHeight = 409.5
$ws.Rows(5).RowHeight = 409.5

# --- Column widths (minor autofit-style refinements) ---
$ws.Columns(1).ColumnWidth = 10.90625
$ws.Columns(2).ColumnWidth = 87.6328125
$ws.Columns(3).ColumnWidth = 75.81640625
$ws.Columns(4).ColumnWidth = 70.1796875
$ws.Columns(5).ColumnWidth = 28.54296875

# --- View state: zoom + active selection land on the new row ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 69
$ws.Range("B5").Select()
